$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Acknowledgement"
$ws.Range("G15").Select()
